$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D (Price) to text format so that numeric-looking
# strings (e.g. "303.54") are not auto-converted to numbers by Excel, matching
# the inline-string cells in the source file.
$priceRange = $ws.Range("D2:D49")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '41.789.36'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').Value = '2.271.04'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '303.54'
$ws.Range('E5').Value = '  +0.33%  '
$ws.Range('D6').Value = '92.56'
$ws.Range('E6').Value = '  +1.17%  '
$ws.Range('D7').Value = '0.530'
$ws.Range('E7').Value = '  +1.83%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').Value = '32.58'
$ws.Range('E10').Value = '  +1.99%  '
$ws.Range('D11').Value = '53.34'
$ws.Range('E11').Value = '  -0.59%  '
$ws.Range('E12').Value = '  +0.34%  '
$ws.Range('E13').Value = '  -1.37%  '
$ws.Range('D14').Value = '6.69'
$ws.Range('E14').Value = '  +1.38%  '
$ws.Range('D15').Value = '2.621.74'
$ws.Range('D16').Value = '14.29'
$ws.Range('E16').Value = '  +1.17%  '
$ws.Range('D17').Value = '2.291.26'
$ws.Range('E17').Value = '  +1.92%  '
$ws.Range('D18').Value = '0.778'
$ws.Range('E18').Value = '  +4.08%  '
$ws.Range('D19').Value = '41.698.63'
$ws.Range('E19').Value = '  +1.41%  '
$ws.Range('D20').Value = '12.45'
$ws.Range('E20').Value = '  +2.42%  '
$ws.Range('D21').Value = '0.0₃0905'
$ws.Range('E21').Value = '  +0.29%  '
$ws.Range('D22').Value = '5.95'
$ws.Range('E22').Value = '  +1.39%  '
$ws.Range('D23').Value = '67.15'
$ws.Range('E23').Value = '  +0.69%  '
$ws.Range('D24').Value = '240.03'
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').Value = '2.58'
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('E26').Value = '  +4.16%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('E28').Value = '  +1.07%  '
$ws.Range('E29').Value = '  -0.74%  '
$ws.Range('D30').Value = '35.59'
$ws.Range('E30').Value = '  +6.40%  '
$ws.Range('D31').Value = '2.06'
$ws.Range('E31').Value = '  -1.56%  '
$ws.Range('D32').Value = '160.80'
$ws.Range('E32').Value = '  +1.60%  '
$ws.Range('E33').Value = '  +1.42%  '
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('E35').Value = '  +1.58%  '
$ws.Range('E36').Value = '  -0.69%  '
$ws.Range('D37').Value = '16.86'
$ws.Range('E37').Value = '  +0.54%  '
$ws.Range('E38').Value = '  +0.50%  '
$ws.Range('E39').Value = '  +1.89%  '
$ws.Range('E40').Value = '  +0.88%  '
$ws.Range('E41').Value = '  +0.87%  '
$ws.Range('D42').Value = '3.92'
$ws.Range('E42').Value = '  -0.35%  '
$ws.Range('D43').Value = '2.002.44'
$ws.Range('E43').Value = '  -3.04%  '
$ws.Range('D44').Value = '19.31'
$ws.Range('E44').Value = '  -4.86%  '
$ws.Range('E45').Value = '  +2.01%  '
$ws.Range('D46').Value = '10.35'
$ws.Range('E46').Value = '  +1.21%  '
$ws.Range('D47').Value = '2.13'
$ws.Range('E47').Value = '  +4.18%  '
$ws.Range('D48').Value = '2.91'
$ws.Range('E48').Value = '  -1.32%  '
$ws.Range('D49').Value = '52.64'
$ws.Range('E49').Value = '  +3.41%  '
$ws.Range('E50').Value = '  +1.00%  '
$ws.Range('E51').Value = '  +0.71%  '

# Restore the default (Normal) style on column D so no stray text-format
# styling is left behind.
$priceRange.Style = "Normal"

Write-Output "done"
